$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instructions")

# Add new note cell in C1 with red font, matching the template update note
$ws.Range("C1").Value = "Template updated 1/17/23"
$ws.Range("C1").Font.Color = 255

# Update the recorded selection on the Instructions sheet (bottomRight pane)
[void]$ws.Range("B15").Select()

# Restore the active sheet back to Meta (the originally active tab)
$wb.Worksheets.Item("Meta").Activate()
